$wb = $excel.ActiveWorkbook

$wsPlayoff = $wb.Worksheets.Item("Playoff Odds")
$wsPlayoff.Range("B2").Value = 48.4
$wsPlayoff.Range("C2").Value = 20.8
$wsPlayoff.Range("E2").Value = 7.6
$wsPlayoff.Range("F2").Value = 4.1
$wsPlayoff.Range("G2").Value = 2.3
$wsPlayoff.Range("H2").Value = 2
$wsPlayoff.Range("I2").Value = 0.6
$wsPlayoff.Range("L2").Value = 97.3
$wsPlayoff.Range("B3").Value = 17.6
$wsPlayoff.Range("C3").Value = 23.6
$wsPlayoff.Range("E3").Value = 12.6
$wsPlayoff.Range("F3").Value = 10.3
$wsPlayoff.Range("G3").Value = 8.8
$wsPlayoff.Range("H3").Value = 5.1
$wsPlayoff.Range("I3").Value = 3.5
$wsPlayoff.Range("J3").Value = 1.1
$wsPlayoff.Range("K3").Value = 0.3
$wsPlayoff.Range("L3").Value = 90
$wsPlayoff.Range("B4").Value = 14.2
$wsPlayoff.Range("C4").Value = 19.3
$wsPlayoff.Range("D4").Value = 17.3
$wsPlayoff.Range("E4").Value = 15.7
$wsPlayoff.Range("F4").Value = 11.4
$wsPlayoff.Range("G4").Value = 9.3
$wsPlayoff.Range("H4").Value = 7.1
$wsPlayoff.Range("I4").Value = 3.5
$wsPlayoff.Range("J4").Value = 1.7
$wsPlayoff.Range("K4").Value = 0.5
$wsPlayoff.Range("L4").Value = 87.2
$wsPlayoff.Range("B5").Value = 7.3
$wsPlayoff.Range("C5").Value = 10.3
$wsPlayoff.Range("D5").Value = 14.7
$wsPlayoff.Range("E5").Value = 15.5
$wsPlayoff.Range("F5").Value = 15.8
$wsPlayoff.Range("G5").Value = 12.9
$wsPlayoff.Range("H5").Value = 13
$wsPlayoff.Range("I5").Value = 6.8
$wsPlayoff.Range("J5").Value = 2.7
$wsPlayoff.Range("K5").Value = 1
$wsPlayoff.Range("L5").Value = 76.5
$wsPlayoff.Range("B6").Value = 4.9
$wsPlayoff.Range("C6").Value = 11.5
$wsPlayoff.Range("D6").Value = 11.9
$wsPlayoff.Range("E6").Value = 12
$wsPlayoff.Range("F6").Value = 13.9
$wsPlayoff.Range("G6").Value = 12.5
$wsPlayoff.Range("H6").Value = 12.5
$wsPlayoff.Range("I6").Value = 11
$wsPlayoff.Range("L6").Value = 66.7
$wsPlayoff.Range("B7").Value = 4.8
$wsPlayoff.Range("C7").Value = 7
$wsPlayoff.Range("D7").Value = 10.9
$wsPlayoff.Range("E7").Value = 14.3
$wsPlayoff.Range("F7").Value = 15.5
$wsPlayoff.Range("G7").Value = 13.9
$wsPlayoff.Range("H7").Value = 13.7
$wsPlayoff.Range("I7").Value = 11.4
$wsPlayoff.Range("J7").Value = 5.8
$wsPlayoff.Range("K7").Value = 2.7
$wsPlayoff.Range("L7").Value = 66.4
$wsPlayoff.Range("C8").Value = 4.3
$wsPlayoff.Range("D8").Value = 8.7
$wsPlayoff.Range("E8").Value = 13.3
$wsPlayoff.Range("F8").Value = 13.8
$wsPlayoff.Range("G8").Value = 17.9
$wsPlayoff.Range("H8").Value = 14.4
$wsPlayoff.Range("I8").Value = 15.1
$wsPlayoff.Range("J8").Value = 7.7
$wsPlayoff.Range("K8").Value = 3
$wsPlayoff.Range("L8").Value = 59.8
$wsPlayoff.Range("C9").Value = 2.8
$wsPlayoff.Range("D9").Value = 4.7
$wsPlayoff.Range("E9").Value = 6.6
$wsPlayoff.Range("F9").Value = 11.8
$wsPlayoff.Range("G9").Value = 14.8
$wsPlayoff.Range("H9").Value = 18.5
$wsPlayoff.Range("I9").Value = 19
$wsPlayoff.Range("J9").Value = 12.8
$wsPlayoff.Range("K9").Value = 8.1
$wsPlayoff.Range("L9").Value = 41.6
$wsPlayoff.Range("C10").Value = 0.3
$wsPlayoff.Range("D10").Value = 0.3
$wsPlayoff.Range("E10").Value = 1.5
$wsPlayoff.Range("F10").Value = 1.6
$wsPlayoff.Range("G10").Value = 4.4
$wsPlayoff.Range("H10").Value = 6.7
$wsPlayoff.Range("I10").Value = 13.5
$wsPlayoff.Range("J10").Value = 30.5
$wsPlayoff.Range("K10").Value = 41.2
$wsPlayoff.Range("L10").Value = 8.1
$wsPlayoff.Range("B11").Value = 0.1
$wsPlayoff.Range("C11").Value = 0.1
$wsPlayoff.Range("E11").Value = 0.9
$wsPlayoff.Range("F11").Value = 1.8
$wsPlayoff.Range("H11").Value = 7
$wsPlayoff.Range("I11").Value = 15.6
$wsPlayoff.Range("J11").Value = 31.4
$wsPlayoff.Range("K11").Value = 39.6
$wsPlayoff.Range("L11").Value = 6.4

$wsRecord = $wb.Worksheets.Item("Record Odds")
$wsRecord.Range("F2").Value = 97.3
$wsRecord.Range("G2").Value = "9.8-4.1-0.1"
$wsRecord.Range("F3").Value = 90
$wsRecord.Range("G3").Value = "8.2-5.7-0.1"
$wsRecord.Range("F4").Value = 87.2
$wsRecord.Range("G4").Value = "8.3-5.7-0.1"
$wsRecord.Range("H4").Value = "9-5"
$wsRecord.Range("F5").Value = 76.5
$wsRecord.Range("G5").Value = "7.7-6.2-0.1"
$wsRecord.Range("F6").Value = 66.7
$wsRecord.Range("G6").Value = "7.2-6.8-0.1"
$wsRecord.Range("F7").Value = 66.4
$wsRecord.Range("G7").Value = "7.1-6.9-0.1"
$wsRecord.Range("F8").Value = 59.8
$wsRecord.Range("G8").Value = "6.5-7.3-0.1"
$wsRecord.Range("F9").Value = 41.6
$wsRecord.Range("G9").Value = "6.0-7.9-0.1"
$wsRecord.Range("H9").Value = "6-8"
$wsRecord.Range("F10").Value = 8.1
$wsRecord.Range("G10").Value = "4.4-9.6-0.1"
$wsRecord.Range("F11").Value = 6.4
$wsRecord.Range("G11").Value = "4.4-9.5-0.1"
$wsRecord.Range("H11").Value = "5-9"

Write-Host "Applied all cell updates"